$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "sahilharpal1234@gmail.com"
$ws.Range("B3").Value = "Sahil Harpal"
$ws.Range("C3").Value = "CA245368"

# Contact number must remain text (the shared string was a digit-only string
# in the source file), so use Excel's leading-apostrophe text marker instead
# of Value = "..." (which would auto-coerce a pure-digit string to a number).
$ws.Range("D3").Value = "'7276801998"

$ws.Range("E3").Value = "Indian Institute of Technology Jodhpur"
